$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 0.2162966666666667
$ws.Range("H2").Value = 0.6488900000000001
$ws.Range("I2").Value = 0.02888548604596741
$ws.Range("J2").Value = 0.0288854860459674
$ws.Range("M2").Value = 0.8366046666666667
$ws.Range("N2").Value = 2.509814
$ws.Range("O2").Value = 0.08025679986157715
$ws.Range("P2").Value = 0.08025679986157715
$ws.Range("Q2").Value = 0.1809548007177778
$ws.Range("R2").Value = 1.62859320646
$ws.Range("S2").Value = 0.002318256672495586
$ws.Range("T2").Value = 0.002318256672495585
$ws.Range("G3").Value = 0.2162966666666667
$ws.Range("H3").Value = 0.6488900000000001
$ws.Range("I3").Value = 0.02888548604596741
$ws.Range("J3").Value = 0.0288854860459674
$ws.Range("M3").Value = 7.939250333333333
$ws.Range("O3").Value = 0.7616247559221037
$ws.Range("P3").Value = 0.7616247559221038
$ws.Range("Q3").Value = 1.717233382932222
$ws.Range("R3").Value = 15.45510044639
$ws.Range("S3").Value = 0.02199990125945126
$ws.Range("T3").Value = 0.02199990125945126
$ws.Range("G4").Value = 0.2162966666666667
$ws.Range("H4").Value = 0.6488900000000001
$ws.Range("I4").Value = 0.02888548604596741
$ws.Range("J4").Value = 0.0288854860459674
$ws.Range("M4").Value = 1.648242
$ws.Range("N4").Value = 4.944726
$ws.Range("O4").Value = 0.1581184442163192
$ws.Range("P4").Value = 0.1581184442163192
$ws.Range("Q4").Value = 0.35650925046
$ws.Range("R4").Value = 3.208583254140001
$ws.Range("S4").Value = 0.004567328114020564
$ws.Range("T4").Value = 0.004567328114020564
$ws.Range("I5").Value = 0.7789723686414617
$ws.Range("J5").Value = 0.7789723686414615
$ws.Range("M5").Value = 0.8366046666666667
$ws.Range("N5").Value = 2.509814
$ws.Range("O5").Value = 0.08025679986157715
$ws.Range("P5").Value = 0.08025679986157715
$ws.Range("Q5").Value = 4.879917530480667
$ws.Range("R5").Value = 43.919257774326
$ws.Range("S5").Value = 0.06251782948775648
$ws.Range("T5").Value = 0.06251782948775647
$ws.Range("I6").Value = 0.7789723686414617
$ws.Range("J6").Value = 0.7789723686414615
$ws.Range("M6").Value = 7.939250333333333
$ws.Range("O6").Value = 0.7616247559221037
$ws.Range("P6").Value = 0.7616247559221038
$ws.Range("S6").Value = 0.5932846401366162
$ws.Range("T6").Value = 0.5932846401366161
$ws.Range("I7").Value = 0.7789723686414617
$ws.Range("J7").Value = 0.7789723686414615
$ws.Range("M7").Value = 1.648242
$ws.Range("N7").Value = 4.944726
$ws.Range("O7").Value = 0.1581184442163192
$ws.Range("P7").Value = 0.1581184442163192
$ws.Range("Q7").Value = 9.614200530726
$ws.Range("R7").Value = 86.52780477653401
$ws.Range("S7").Value = 0.123169899017089
$ws.Range("T7").Value = 0.123169899017089
$ws.Range("G8").Value = 1.438774666666667
$ws.Range("H8").Value = 4.316324
$ws.Range("I8").Value = 0.192142145312571
$ws.Range("J8").Value = 0.192142145312571
$ws.Range("M8").Value = 0.8366046666666667
$ws.Range("N8").Value = 2.509814
$ws.Range("O8").Value = 0.08025679986157715
$ws.Range("P8").Value = 0.08025679986157715
$ws.Range("Q8").Value = 1.203685600415111
$ws.Range("R8").Value = 10.833170403736
$ws.Range("S8").Value = 0.01542071370132508
$ws.Range("T8").Value = 0.01542071370132508
$ws.Range("G9").Value = 1.438774666666667
$ws.Range("H9").Value = 4.316324
$ws.Range("I9").Value = 0.192142145312571
$ws.Range("J9").Value = 0.192142145312571
$ws.Range("M9").Value = 7.939250333333333
$ws.Range("O9").Value = 0.7616247559221037
$ws.Range("P9").Value = 0.7616247559221038
$ws.Range("Q9").Value = 11.42279225192489
$ws.Range("R9").Value = 102.805130267324
$ws.Range("S9").Value = 0.1463402145260362
$ws.Range("T9").Value = 0.1463402145260363
$ws.Range("G10").Value = 1.438774666666667
$ws.Range("H10").Value = 4.316324
$ws.Range("I10").Value = 0.192142145312571
$ws.Range("J10").Value = 0.192142145312571
$ws.Range("M10").Value = 1.648242
$ws.Range("N10").Value = 4.944726
$ws.Range("O10").Value = 0.1581184442163192
$ws.Range("P10").Value = 0.1581184442163192
$ws.Range("Q10").Value = 2.371448834136
$ws.Range("R10").Value = 21.343039507224
$ws.Range("S10").Value = 0.03038121708520965
$ws.Range("T10").Value = 0.03038121708520965
